$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("C2").Value = 26.5
$ws.Range("D2").Value = 105416832
$ws.Range("E2").Value = 103889024
$ws.Range("F2").Value = 1527808
$ws.Range("G2").Value = 1.45

# Add new rows 3-8 with modulation names (column A) and status (column H)
$ws.Range("A3").Value = "QAM64 2/3"
$ws.Range("H3").Value = "Ошибка поиска модуляции"

$ws.Range("A4").Value = "QAM16 3/4"
$ws.Range("H4").Value = "Ошибка поиска модуляции"

$ws.Range("A5").Value = "QAM16 1/2"
$ws.Range("H5").Value = "Ошибка поиска модуляции"

$ws.Range("A6").Value = "QPSK 3/4"
$ws.Range("H6").Value = "Ошибка поиска модуляции"

$ws.Range("A7").Value = "QPSK 1/2"
$ws.Range("H7").Value = "Ошибка поиска модуляции"

$ws.Range("A8").Value = "BPSK 1/2"
$ws.Range("H8").Value = "Ошибка поиска модуляции"
